# Applies the CircadiPy re-run (cosinor_2_sawtooth_10) value refresh to Sheet1.
# PowerShell's literal parser chokes on bare scientific-notation tokens like
# "5.062002705891189e-09" (the tokenizer reads the "e-09" part as a separate
# bareword), so every such constant is built via [System.Convert]::ToDouble on
# a string literal instead of being typed as a numeric literal directly.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# ---- Row 2 ----
$ws.Range("E2").Value = 22.49000000000008
$ws.Range("G2").Value = [System.Convert]::ToDouble("5.062002705891189e-09")
$ws.Range("H2").Value = [System.Convert]::ToDouble("6.065825708390093e-08")
$ws.Range("K2").Value = 5.897127839255191
$ws.Range("L2").Value = "[3.487808993591784, 8.306446684918598]"
$ws.Range("M2").Value = [System.Convert]::ToDouble("2.084699568527526e-06")
$ws.Range("N2").Value = [System.Convert]::ToDouble("3.750292164816926e-06")
$ws.Range("O2").Value = -1.119526511189155
$ws.Range("P2").Value = "[-1.547210796362541, -0.69184222601577]"
$ws.Range("Q2").Value = [System.Convert]::ToDouble("4.07861643791918e-07")
$ws.Range("R2").Value = [System.Convert]::ToDouble("8.157232875838361e-07")
$ws.Range("S2").Value = 10.58006233290856
$ws.Range("T2").Value = "[9.285865315435952, 11.874259350381177]"
$ws.Range("W2").Value = 4.00722722722724
$ws.Range("X2").Value = 2.476376376376383
$ws.Range("Y2").Value = 5.538078078078097

# ---- Row 3 ----
$ws.Range("E3").Value = 23.2900000000002
$ws.Range("G3").Value = [System.Convert]::ToDouble("4.312395995853535e-09")
$ws.Range("H3").Value = [System.Convert]::ToDouble("6.065825708390093e-08")
$ws.Range("I3").ClearContents()
$ws.Range("K3").Value = 6.519191894863027
$ws.Range("L3").Value = "[3.7928016248947163, 9.245582164831337]"
$ws.Range("M3").Value = [System.Convert]::ToDouble("3.750292164816926e-06")
$ws.Range("N3").Value = [System.Convert]::ToDouble("3.750292164816926e-06")
$ws.Range("O3").Value = 0.9874475407679633
$ws.Range("P3").Value = "[0.5723422051585008, 1.4025528763774258]"
$ws.Range("Q3").Value = [System.Convert]::ToDouble("4.192374915001906e-06")
$ws.Range("R3").Value = [System.Convert]::ToDouble("4.192374915001906e-06")
$ws.Range("S3").Value = 10.24158687569598
$ws.Range("T3").Value = "[8.80943303589412, 11.67374071549784]"
$ws.Range("W3").Value = 19.62980980980998
$ws.Range("X3").Value = 18.09113113113128
$ws.Range("Y3").Value = 21.16848848848867
